$d = $word.ActiveDocument

# 1) "O Sumário Executivo aborda, do Plano de Negócios de maneira sucinta..."
#    -> "O Sumário Executivo aborda, de maneira sucinta..."
$d.Content.Find.Execute("do Plano de Negócios de maneira sucinta", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "de maneira sucinta", 2)

# 2) "...no modelo dissertativa argumentativa, onde estas serão avaliadas por professores. Estre projeto..."
#    -> "...no modelo dissertativo argumentativo, sendo que estas serão avaliadas por professores. Este projeto..."
$d.Content.Find.Execute("dissertativa argumentativa, onde estas serão avaliadas por professores. Estre projeto", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "dissertativo argumentativo, sendo que estas serão avaliadas por professores. Este projeto", 2)

# 3) Rewrite the "Missão e Valores do Projeto" explanatory sentence
$d.Content.Find.Execute( `
  "colaborar para o melhoramento de praticantes da redação no modelo ENEM, visando o aprendizado destes praticantes. Os desenvolvedores do projeto acreditam que a educação é o princípio para uma sociedade próspera, com tal filosofia em mente foi escolhido, para Trabalho de Conclusão de Curso do curso Técnico de Informática do SENAI, algo relacionado à educação.", `
  $false, $false, $false, $false, $false, `
  $true, 1, $false, `
  "colaborar para o melhoramento da prática de desenvolvimento de redações no modelo ENEM, visando o aprendizado dos praticantes. Os desenvolvedores do projeto acreditam que a educação é o princípio para uma sociedade próspera e, com tal filosofia em mente, foi escolhido, para Trabalho de Conclusão de Curso do curso Técnico de Informática do SENAI - Portão, algo relacionado à educação.", 2)

# 4) " de mercado em que abrigará o projeto" -> " de mercado que abrigará o projeto"
$d.Content.Find.Execute("de mercado em que abrigará o projeto", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "de mercado que abrigará o projeto", 2)

# 5) Move the "_GoBack" bookmark from the end of the "Analise de Mercado" paragraph
#    to between "c" and "urso" in "curso Técnico de Informática do SENAI - Portão"
#    (re-adding a bookmark with the same name relocates it).
$rng = $d.Content
$rng.Find.Execute("curso Técnico de Informática do SENAI - Portão", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$bm = $d.Range($rng.Start + 1, $rng.Start + 1)
$d.Bookmarks.Add("_GoBack", $bm)

Write-Output "done"
